$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pass 1: touch the handful of brand-new label strings first, in the same
# order the target workbook introduces them, so the saved sharedStrings table
# allocates identical indices/order (cosmetic, but keeps the diff minimal). ---
$ws.Range("A81").Value = "proximity.1.above1"
$ws.Range("B80").Value = "anteater_A"
$ws.Range("C80").Value = "anteater_B"
$ws.Range("A83").Value = "proximity.1.below1"
$ws.Range("A90").Value = "proximity.2.above1"
$ws.Range("B91").Value = "NONE"
$ws.Range("A92").Value = "proximity.2.below1"
$ws.Range("A97").Value = "proximity.3.above1"
$ws.Range("A99").Value = "proximity.3.below1"

# --- Pass 2: seed three new cell styles in off-sheet helper cells, in the
# same order they need to appear in the saved cellXfs table:
#   18 = default font, horizontal-left alignment
#   19 = default font, alignment flag applied but no explicit sub-properties
#   20 = light-blue font (matches existing fontId 5), horizontal-left alignment
# then Copy / PasteSpecial(xlPasteFormats) them onto the real target cells so
# the column-level style (style 1, centered) that column A otherwise inherits
# does not bleed through. ---
$xlPasteFormats = -4122
$xlLeft = -4131

$ws.Range("Z200").Value = "tmp18"
$ws.Range("Z200").HorizontalAlignment = $xlLeft

$ws.Range("Z201").Value = "tmp19"
$ws.Range("Z201").WrapText = $true
$ws.Range("Z201").WrapText = $false

$ws.Range("Z202").Value = "tmp20"
$ws.Range("Z202").HorizontalAlignment = $xlLeft
$ws.Range("Z202").Font.Color = 15773696

$style18 = $ws.Range("Z200")
$style19 = $ws.Range("Z201")
$style20 = $ws.Range("Z202")
$lightBlue = 15773696

# --- Pass 3: fill in every cells final value + style for the new proximity-pair block. ---
$ws.Range("B80").Value = "anteater_A"
$ws.Range("C80").Value = "anteater_B"

$style19.Copy()
$ws.Range("A81").PasteSpecial($xlPasteFormats)
$ws.Range("A81").Value = "proximity.1.above1"

$style19.Copy()
$ws.Range("A82").PasteSpecial($xlPasteFormats)
$style20.Copy()
$ws.Range("B82").PasteSpecial($xlPasteFormats)
$ws.Range("B82").Value = "Kyle"
$style18.Copy()
$ws.Range("C82").PasteSpecial($xlPasteFormats)
$ws.Range("C82").Value = "Christoffer"

$style19.Copy()
$ws.Range("A83").PasteSpecial($xlPasteFormats)
$ws.Range("A83").Value = "proximity.1.below1"

$style19.Copy()
$ws.Range("A84").PasteSpecial($xlPasteFormats)
$ws.Range("B84").Value = "Elaine"
$style18.Copy()
$ws.Range("C84").PasteSpecial($xlPasteFormats)
$ws.Range("C84").Value = "Christoffer"

$style19.Copy()
$ws.Range("A85").PasteSpecial($xlPasteFormats)
$ws.Range("B85").Value = "Kyle"
$ws.Range("B85").Font.Color = $lightBlue
$style18.Copy()
$ws.Range("C85").PasteSpecial($xlPasteFormats)
$ws.Range("C85").Value = "Bumpus"

$style19.Copy()
$ws.Range("A86").PasteSpecial($xlPasteFormats)
$ws.Range("B86").Value = "Little Rick"
$ws.Range("B86").Font.Color = $lightBlue
$style18.Copy()
$ws.Range("C86").PasteSpecial($xlPasteFormats)
$ws.Range("C86").Value = "Elaine"

$style19.Copy()
$ws.Range("A87").PasteSpecial($xlPasteFormats)
$ws.Range("B87").Value = "Makao"
$style18.Copy()
$ws.Range("C87").PasteSpecial($xlPasteFormats)
$ws.Range("C87").Value = "Bumpus"

$style19.Copy()
$ws.Range("A88").PasteSpecial($xlPasteFormats)
$ws.Range("B88").Value = "Puji"
$style18.Copy()
$ws.Range("C88").PasteSpecial($xlPasteFormats)
$ws.Range("C88").Value = "Bumpus"

$style19.Copy()
$ws.Range("A89").PasteSpecial($xlPasteFormats)
$ws.Range("B89").Value = "Rodolfo"
$style18.Copy()
$ws.Range("C89").PasteSpecial($xlPasteFormats)
$ws.Range("C89").Value = "Elaine"

$style19.Copy()
$ws.Range("A90").PasteSpecial($xlPasteFormats)
$ws.Range("A90").Value = "proximity.2.above1"

$style19.Copy()
$ws.Range("A91").PasteSpecial($xlPasteFormats)
$ws.Range("B91").Value = "NONE"

$style18.Copy()
$ws.Range("A92").PasteSpecial($xlPasteFormats)
$ws.Range("A92").Value = "proximity.2.below1"

$style18.Copy()
$ws.Range("A93").PasteSpecial($xlPasteFormats)
$ws.Range("B93").Value = "Larry"
$style18.Copy()
$ws.Range("C93").PasteSpecial($xlPasteFormats)
$ws.Range("C93").Value = "Annie"

$style18.Copy()
$ws.Range("A94").PasteSpecial($xlPasteFormats)
$ws.Range("B94").Value = "Reid"
$ws.Range("B94").Font.Color = $lightBlue
$style18.Copy()
$ws.Range("C94").PasteSpecial($xlPasteFormats)
$ws.Range("C94").Value = "Larry"

$style18.Copy()
$ws.Range("A95").PasteSpecial($xlPasteFormats)
$ws.Range("B95").Value = "Thomas"
$style18.Copy()
$ws.Range("C95").PasteSpecial($xlPasteFormats)
$ws.Range("C95").Value = "Margaret"

$style18.Copy()
$ws.Range("A96").PasteSpecial($xlPasteFormats)
$ws.Range("B96").Value = "Thomas"
$style20.Copy()
$ws.Range("C96").PasteSpecial($xlPasteFormats)
$ws.Range("C96").Value = "Reid"

$style18.Copy()
$ws.Range("A97").PasteSpecial($xlPasteFormats)
$ws.Range("A97").Value = "proximity.3.above1"

$style18.Copy()
$ws.Range("A98").PasteSpecial($xlPasteFormats)
$ws.Range("B98").Value = "NONE"

$style18.Copy()
$ws.Range("A99").PasteSpecial($xlPasteFormats)
$ws.Range("A99").Value = "proximity.3.below1"

$style18.Copy()
$ws.Range("A100").PasteSpecial($xlPasteFormats)
$ws.Range("B100").Value = "Sheron"
$style18.Copy()
$ws.Range("C100").PasteSpecial($xlPasteFormats)
$ws.Range("C100").Value = "Maria"

$style18.Copy()
$ws.Range("A101").PasteSpecial($xlPasteFormats)

# clean up helper cells + clipboard marching-ants state
$ws.Range("Z200:Z202").Clear()
$excel.CutCopyMode = $false
